$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("C15").Value = 7.94

# Row 16
$ws.Range("C16").Value = 7.94

# Row 17
$ws.Range("C17").Value = 7.94

# Row 18
$ws.Range("C18").Value = 7.94

# Row 19
$ws.Range("C19").Value = 7.94

# Row 20
$ws.Range("C20").Value = 7.94
$ws.Range("E20").Value = 13.7
$ws.Range("F20").Value = 7.9
$ws.Range("H20").Value = 10.894

# Row 21
$ws.Range("C21").Value = 7.94

# Row 22
$ws.Range("C22").Value = 7.94

# Row 23
$ws.Range("C23").Value = 7.94

# Row 24
$ws.Range("C24").Value = 7.94
$ws.Range("K24").Value = 94.6
$ws.Range("N24").Value = 4.8471000000000002
$ws.Range("Q24").Value = 11.5
$ws.Range("T24").Value = 19.623000000000001

# Row 25
$ws.Range("C25").Value = 7.94
$ws.Range("E25").Value = 13.6
$ws.Range("F25").Value = 10.8
$ws.Range("H25").Value = 10.896000000000001
$ws.Range("K25").Value = 129.19999999999999
$ws.Range("L25").Value = 18.2
$ws.Range("N25").Value = 4.8471000000000002
$ws.Range("Q25").Value = 55.2
$ws.Range("R25").Value = 10.4
$ws.Range("T25").Value = 19.623000000000001
$ws.Range("W25").Value = 253.9
$ws.Range("X25").Value = 51
$ws.Range("Z25").Value = 12.440899999999999

# Row 26
$ws.Range("C26").Value = 7.94
$ws.Range("E26").Value = 14.2
$ws.Range("H26").Value = 10.896000000000001
$ws.Range("K26").Value = 90.8
$ws.Range("L26").Value = 19.8
$ws.Range("N26").Value = 4.8471000000000002
$ws.Range("Q26").Value = 43.5
$ws.Range("R26").Value = 10.6
$ws.Range("T26").Value = 19.619
$ws.Range("W26").Value = 139.19999999999999
$ws.Range("X26").Value = 48
$ws.Range("Z26").Value = 12.4429

# Row 27
$ws.Range("C27").Value = 7.94
$ws.Range("E27").Value = 12.6
$ws.Range("F27").Value = 6.9
$ws.Range("H27").Value = 10.894
$ws.Range("K27").Value = 88.5
$ws.Range("L27").Value = 19.3
$ws.Range("N27").Value = 4.8471000000000002

# Update selection to match the diff
$ws.Range("E31").Select()
